$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.420.26"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "3.378.28"
$ws.Range("E3").Value = "  +0.34%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "575.22"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("D6").Value = "136.75"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("D8").Value = "3.378.38"
$ws.Range("E8").Value = "  +0.33%  "

$ws.Range("E9").Value = "  -0.45%  "

$ws.Range("D10").Value = "7.48"
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("E11").Value = "  +1.95%  "

$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").Value = "3.957.61"
$ws.Range("E13").Value = "  +0.51%  "

$ws.Range("E14").Value = "  +2.50%  "

$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  +2.14%  "

$ws.Range("D16").Value = "26.01"
$ws.Range("E16").Value = "  +3.59%  "

$ws.Range("D17").Value = "3.380.31"
$ws.Range("E17").Value = "  +0.62%  "

$ws.Range("D18").Value = "61.633.32"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("D19").Value = "14.08"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("D21").Value = "9.37"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").Value = "376.66"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("D23").Value = "0.556"
$ws.Range("E23").Value = "  -2.84%  "

$ws.Range("D24").Value = "3.518.63"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("E25").Value = "  -0.24%  "

$ws.Range("E26").Value = "  +7.23%  "

$ws.Range("D27").Value = "71.37"
$ws.Range("E27").Value = "  +0.76%  "

$ws.Range("D28").Value = "1.74"
$ws.Range("E28").Value = "  +6.19%  "

$ws.Range("D29").Value = "7.48"
$ws.Range("E29").Value = "  -3.25%  "

$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("D31").Value = "8.24"
$ws.Range("E31").Value = "  +1.97%  "

$ws.Range("D32").Value = "'0.160"
$ws.Range("E32").Value = "  +3.74%  "

$ws.Range("E33").Value = "  +1.79%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").Value = "23.45"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "5.27"
$ws.Range("E36").Value = "  -5.07%  "

$ws.Range("E37").Value = "  -1.28%  "

$ws.Range("E38").Value = "  -0.75%  "

$ws.Range("D39").Value = "164.84"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("D40").Value = "0.0774"
$ws.Range("E40").Value = "  -2.24%  "

$ws.Range("E41").Value = "  +0.15%  "

$ws.Range("D42").Value = "0.774"
$ws.Range("E42").Value = "  +2.22%  "

$ws.Range("D43").Value = "1.71"
$ws.Range("E43").Value = "  +6.37%  "

$ws.Range("E44").Value = "  +1.43%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "41.51"
$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'4.40"
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("D47").Value = "24.74"
$ws.Range("E47").Value = "  +7.89%  "

$ws.Range("E48").Value = "  -1.72%  "

$ws.Range("D49").Value = "22.73"
$ws.Range("E49").Value = "  -1.39%  "

$ws.Range("D50").Value = "2.346.34"
$ws.Range("E50").Value = "  +4.38%  "

$ws.Range("E51").Value = "  +1.55%  "
